$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "330.33"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "7.22%"
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "19"
$ws.Cells.Item(2, 7).Style = "Normal"

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "40.73"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "12.78%"
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "19"
$ws.Cells.Item(3, 7).Style = "Normal"

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "6.117"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "19.74%"
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = "19"
$ws.Cells.Item(4, 7).Style = "Normal"

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.08198"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "6.24%"
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(5, 7).NumberFormat = "@"
$ws.Cells.Item(5, 7).Value = "19"
$ws.Cells.Item(5, 7).Style = "Normal"

# Row 6
$ws.Cells.Item(6, 2).Value = "GateToken"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "4.587"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "4.54%"
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(6, 7).NumberFormat = "@"
$ws.Cells.Item(6, 7).Value = "19"
$ws.Cells.Item(6, 7).Style = "Normal"

# Row 7
$ws.Cells.Item(7, 2).Value = "KuCoinToken"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "8.807"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "6.18%"
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Cells.Item(7, 7).Value = "19"
$ws.Cells.Item(7, 7).Style = "Normal"

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "1.974"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "7.24%"
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = "19"
$ws.Cells.Item(8, 7).Style = "Normal"

# Row 9
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = "19"
$ws.Cells.Item(9, 7).Style = "Normal"

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.9515"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "3.15%"
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(10, 7).NumberFormat = "@"
$ws.Cells.Item(10, 7).Value = "19"
$ws.Cells.Item(10, 7).Style = "Normal"

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.1358"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "22.11%"
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(11, 7).NumberFormat = "@"
$ws.Cells.Item(11, 7).Value = "19"
$ws.Cells.Item(11, 7).Style = "Normal"

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.2010"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "8.13%"
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(12, 7).NumberFormat = "@"
$ws.Cells.Item(12, 7).Value = "19"
$ws.Cells.Item(12, 7).Style = "Normal"

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "10.46"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "64.71%"
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(13, 7).NumberFormat = "@"
$ws.Cells.Item(13, 7).Value = "19"
$ws.Cells.Item(13, 7).Style = "Normal"

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.09250"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "5.56%"
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(14, 7).NumberFormat = "@"
$ws.Cells.Item(14, 7).Value = "19"
$ws.Cells.Item(14, 7).Style = "Normal"

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.03543"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "6.54%"
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(15, 7).NumberFormat = "@"
$ws.Cells.Item(15, 7).Value = "19"
$ws.Cells.Item(15, 7).Style = "Normal"

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.09641"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "1.14%"
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = "19"
$ws.Cells.Item(16, 7).Style = "Normal"

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.001308"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "-5.53%"
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(17, 7).NumberFormat = "@"
$ws.Cells.Item(17, 7).Value = "19"
$ws.Cells.Item(17, 7).Style = "Normal"

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.006323"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "2.76%"
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(18, 7).NumberFormat = "@"
$ws.Cells.Item(18, 7).Value = "19"
$ws.Cells.Item(18, 7).Style = "Normal"

# Row 19
$ws.Cells.Item(19, 2).Value = "LEO"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "3.358"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "-0.02%"
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(19, 7).NumberFormat = "@"
$ws.Cells.Item(19, 7).Value = "19"
$ws.Cells.Item(19, 7).Style = "Normal"

# Row 20
$ws.Cells.Item(20, 2).Value = "BitpandaEcosystemToken"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.3524"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "2.63%"
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(20, 7).NumberFormat = "@"
$ws.Cells.Item(20, 7).Value = "19"
$ws.Cells.Item(20, 7).Style = "Normal"

# Row 21
$ws.Cells.Item(21, 2).Value = "ProBitToken"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.1435"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "10.02%"
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(21, 7).NumberFormat = "@"
$ws.Cells.Item(21, 7).Value = "19"
$ws.Cells.Item(21, 7).Style = "Normal"

# Row 22
$ws.Cells.Item(22, 2).Value = "ZBToken"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.2432"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "5.16%"
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(22, 7).NumberFormat = "@"
$ws.Cells.Item(22, 7).Value = "19"
$ws.Cells.Item(22, 7).Style = "Normal"

# Row 23
$ws.Cells.Item(23, 2).Value = "CoinExToken"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.04441"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "2.52%"
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(23, 7).NumberFormat = "@"
$ws.Cells.Item(23, 7).Value = "19"
$ws.Cells.Item(23, 7).Style = "Normal"

# Row 24
$ws.Cells.Item(24, 2).Value = "BitKan"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.001258"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "4.75%"
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(24, 7).NumberFormat = "@"
$ws.Cells.Item(24, 7).Value = "19"
$ws.Cells.Item(24, 7).Style = "Normal"

# Row 25
$ws.Cells.Item(25, 2).Value = "HotbitToken"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.004431"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "3.93%"
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(25, 7).NumberFormat = "@"
$ws.Cells.Item(25, 7).Value = "19"
$ws.Cells.Item(25, 7).Style = "Normal"

# Row 26
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "-18.26%"
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(26, 7).NumberFormat = "@"
$ws.Cells.Item(26, 7).Value = "19"
$ws.Cells.Item(26, 7).Style = "Normal"

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.0003977"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = "37.05%"
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(27, 7).NumberFormat = "@"
$ws.Cells.Item(27, 7).Value = "19"
$ws.Cells.Item(27, 7).Style = "Normal"

# Row 28
$ws.Cells.Item(28, 7).NumberFormat = "@"
$ws.Cells.Item(28, 7).Value = "19"
$ws.Cells.Item(28, 7).Style = "Normal"

# Row 29
$ws.Cells.Item(29, 7).NumberFormat = "@"
$ws.Cells.Item(29, 7).Value = "19"
$ws.Cells.Item(29, 7).Style = "Normal"

# Row 30
$ws.Cells.Item(30, 7).NumberFormat = "@"
$ws.Cells.Item(30, 7).Value = "19"
$ws.Cells.Item(30, 7).Style = "Normal"

# Row 31
$ws.Cells.Item(31, 7).NumberFormat = "@"
$ws.Cells.Item(31, 7).Value = "19"
$ws.Cells.Item(31, 7).Style = "Normal"

# Row 32
$ws.Cells.Item(32, 7).NumberFormat = "@"
$ws.Cells.Item(32, 7).Value = "19"
$ws.Cells.Item(32, 7).Style = "Normal"

# Row 33
$ws.Cells.Item(33, 7).NumberFormat = "@"
$ws.Cells.Item(33, 7).Value = "19"
$ws.Cells.Item(33, 7).Style = "Normal"

# Row 34
$ws.Cells.Item(34, 7).NumberFormat = "@"
$ws.Cells.Item(34, 7).Value = "19"
$ws.Cells.Item(34, 7).Style = "Normal"

# Row 35
$ws.Cells.Item(35, 7).NumberFormat = "@"
$ws.Cells.Item(35, 7).Value = "19"
$ws.Cells.Item(35, 7).Style = "Normal"

# Row 36
$ws.Cells.Item(36, 7).NumberFormat = "@"
$ws.Cells.Item(36, 7).Value = "19"
$ws.Cells.Item(36, 7).Style = "Normal"

# Row 37
$ws.Cells.Item(37, 7).NumberFormat = "@"
$ws.Cells.Item(37, 7).Value = "19"
$ws.Cells.Item(37, 7).Style = "Normal"

# Row 38
$ws.Cells.Item(38, 7).NumberFormat = "@"
$ws.Cells.Item(38, 7).Value = "19"
$ws.Cells.Item(38, 7).Style = "Normal"

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.02518"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "20.21%"
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(39, 7).NumberFormat = "@"
$ws.Cells.Item(39, 7).Value = "19"
$ws.Cells.Item(39, 7).Style = "Normal"

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.05295"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "6.94%"
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(40, 7).NumberFormat = "@"
$ws.Cells.Item(40, 7).Value = "19"
$ws.Cells.Item(40, 7).Style = "Normal"

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.007470"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "-0.78%"
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(41, 7).NumberFormat = "@"
$ws.Cells.Item(41, 7).Value = "19"
$ws.Cells.Item(41, 7).Style = "Normal"

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.1454"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "7.57%"
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(42, 7).NumberFormat = "@"
$ws.Cells.Item(42, 7).Value = "19"
$ws.Cells.Item(42, 7).Style = "Normal"

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.008938"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "5.87%"
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(43, 7).NumberFormat = "@"
$ws.Cells.Item(43, 7).Value = "19"
$ws.Cells.Item(43, 7).Style = "Normal"

# Row 44
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "-1.15%"
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(44, 7).NumberFormat = "@"
$ws.Cells.Item(44, 7).Value = "19"
$ws.Cells.Item(44, 7).Style = "Normal"

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.01054"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "24.87%"
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(45, 7).NumberFormat = "@"
$ws.Cells.Item(45, 7).Value = "19"
$ws.Cells.Item(45, 7).Style = "Normal"

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.00006792"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "7.51%"
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(46, 7).NumberFormat = "@"
$ws.Cells.Item(46, 7).Value = "19"
$ws.Cells.Item(46, 7).Style = "Normal"

# Row 47
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "-0.20%"
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(47, 7).NumberFormat = "@"
$ws.Cells.Item(47, 7).Value = "19"
$ws.Cells.Item(47, 7).Style = "Normal"

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.003484"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = "21.44%"
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(48, 7).NumberFormat = "@"
$ws.Cells.Item(48, 7).Value = "19"
$ws.Cells.Item(48, 7).Style = "Normal"

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.001798"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 7).NumberFormat = "@"
$ws.Cells.Item(49, 7).Value = "19"
$ws.Cells.Item(49, 7).Style = "Normal"

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.00002097"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = "-0.20%"
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(50, 7).NumberFormat = "@"
$ws.Cells.Item(50, 7).Value = "19"
$ws.Cells.Item(50, 7).Style = "Normal"

# Row 51
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = "-0.20%"
$ws.Cells.Item(51, 5).Style = "Normal"
$ws.Cells.Item(51, 7).NumberFormat = "@"
$ws.Cells.Item(51, 7).Value = "19"
$ws.Cells.Item(51, 7).Style = "Normal"

